# Add I0 and IF columns (I and J) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells
$ws.Cells.Item(1, 9).Value = "I0"
$ws.Cells.Item(1, 10).Value = "IF"

# Copy the formatting (bold, border, centered) from H1 onto I1/J1
$ws.Cells.Item(1, 8).Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data, encoded as "row,I,J" pairs
$data = "2,7,8;3,6,6;4,7,7;5,4,4;6,5,6;7,8,8;8,5,5;9,7,7;10,7,8;11,6,6;12,7,7;13,6,7;14,6,6;15,6,7;16,7,7;17,7,7;18,7,7;19,7,7;20,6,6;21,7,7;22,9,9;23,8,8;24,7,7;25,6,6;26,6,7;27,7,7;28,11,11;29,8,8;30,8,8;31,7,7;32,7,8;33,5,6;34,7,8;35,6,6;36,8,8;37,6,6;38,6,6;39,6,6;40,7,7;41,7,8;42,6,6;43,6,6;44,8,8;45,7,7;46,5,6;47,7,8;48,6,6;49,7,7;50,7,7;51,4,4;52,7,7;53,3,4;54,6,6;55,7,7;56,6,6;57,5,5;58,9,9;59,8,8;60,6,6;61,5,6;62,6,6;63,6,7;64,6,6;65,8,8;66,8,8;67,8,8;68,8,8;69,7,8;70,8,8;71,6,7;72,8,9;73,6,7;74,5,5;75,6,7;76,6,6;77,7,7;78,6,7;79,7,8;80,6,6;81,6,6;82,6,7;83,9,9;84,6,6;85,8,8;86,4,4;87,7,7;88,9,9;89,7,7;90,6,6;91,5,6;92,8,8;93,6,6;94,3,3;95,3,3;96,3,3"

foreach ($entry in $data.Split(";")) {
    $fields = $entry.Split(",")
    $row = [int]$fields[0]
    $iVal = [int]$fields[1]
    $jVal = [int]$fields[2]
    $ws.Cells.Item($row, 9).Value = $iVal
    $ws.Cells.Item($row, 10).Value = $jVal
}
